$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new columns before column D. This shifts existing D:K data to F:M,
# matching the two newly-reported quarters (2018-12-31, 2018-09-30) being added in D:E.
$ws.Range("D:E").Insert()

# The blank new D:E columns do not inherit number formats/fonts from the (now-shifted)
# data columns, so copy formatting from column F (old column D) across each table block.
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:F35").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)

$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F39:F77").Copy()
$ws.Range("D39:E77").PasteSpecial(-4122)

$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("F81:F102").Copy()
$ws.Range("D81:E102").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# New columns also need a column width; match the (now-shifted) neighbouring data columns
$ws.Range("D1:E1").ColumnWidth = $ws.Range("F1").ColumnWidth

# Populate new columns D and E with the newly-reported quarter data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1219200
$ws.Range("E8").Value = 1035600
$ws.Range("D9").Value = 677300
$ws.Range("E9").Value = 588100
$ws.Range("D10").Value = 541900
$ws.Range("E10").Value = 447500
$ws.Range("D12").Value = 77700
$ws.Range("E12").Value = 70900
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 62800
$ws.Range("E14").Value = 50100
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 1019300
$ws.Range("E17").Value = 895400
$ws.Range("D18").Value = 199900
$ws.Range("E18").Value = 140200
$ws.Range("D20").Value = -4900
$ws.Range("E20").Value = 2800
$ws.Range("D21").Value = 238500
$ws.Range("E21").Value = 185400
$ws.Range("D22").Value = 32800
$ws.Range("E22").Value = 33400
$ws.Range("D23").Value = 162200
$ws.Range("E23").Value = 109600
$ws.Range("D24").Value = -16900
$ws.Range("E24").Value = 12600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 179100
$ws.Range("E26").Value = 97000
$ws.Range("D27").Value = 176100
$ws.Range("E27").Value = 94600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -144500
$ws.Range("E29").Value = -22200
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4900
$ws.Range("E32").Value = -2800
$ws.Range("D33").Value = 31600
$ws.Range("E33").Value = 72400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 31600
$ws.Range("E35").Value = 72400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 161700
$ws.Range("E41").Value = 176500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2285200
$ws.Range("E43").Value = 1901900
$ws.Range("D44").Value = 1097300
$ws.Range("E44").Value = 995400
$ws.Range("D45").Value = 486000
$ws.Range("E45").Value = 417400
$ws.Range("D46").Value = 4030200
$ws.Range("E46").Value = 3491200
$ws.Range("D47").Value = 700
$ws.Range("E47").Value = 81900
$ws.Range("D48").Value = 1032600
$ws.Range("E48").Value = 1002700
$ws.Range("D49").Value = 4172400
$ws.Range("E49").Value = 4206300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 738400
$ws.Range("E52").Value = 624800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 9974300
$ws.Range("E54").Value = 9406900
$ws.Range("D57").Value = 867500
$ws.Range("E57").Value = 784200
$ws.Range("D58").Value = 547700
$ws.Range("E58").Value = 149400
$ws.Range("D59").Value = 1578200
$ws.Range("E59").Value = 1337900
$ws.Range("D60").Value = 2993400
$ws.Range("E60").Value = 2271500
$ws.Range("D61").Value = 2179000
$ws.Range("E61").Value = 2593300
$ws.Range("D62").Value = 1591500
$ws.Range("E62").Value = 1463800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6853200
$ws.Range("E66").Value = 6360500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 4334300
$ws.Range("E72").Value = 4355100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3121100
$ws.Range("E76").Value = 3046400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 31600
$ws.Range("E81").Value = 72400
$ws.Range("D83").Value = 43500
$ws.Range("E83").Value = 42400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -45600
$ws.Range("E89").Value = 220600
$ws.Range("D91").Value = -68000
$ws.Range("E91").Value = -35300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -93400
$ws.Range("E94").Value = -41400
$ws.Range("D96").Value = -22300
$ws.Range("E96").Value = -22300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 122200
$ws.Range("E100").Value = -330300
$ws.Range("D101").Value = 2000
$ws.Range("E101").Value = 1200
$ws.Range("D102").Value = -14800
$ws.Range("E102").Value = -149900

# Apply restated prior-quarter figures (columns H and I) that came with this update
$ws.Range("H14").Value = 102800
$ws.Range("I14").Value = 53200
$ws.Range("H17").Value = 906600
$ws.Range("I17").Value = 572100
$ws.Range("H18").Value = 73000
$ws.Range("I18").Value = 74100
$ws.Range("H20").Value = -29200
$ws.Range("I20").Value = 3600
$ws.Range("H32").Value = 29200
$ws.Range("I32").Value = -3600
